# Fix shear max design, simplify shear design, add cirsoc shear design rebar.
# Updates the numeric results in the "Limit checks", "Shear reinforcement
# strength" and "Shear strength" tables to reflect the corrected calculation.

$d = $word.ActiveDocument

function Replace-Value([string]$old, [string]$new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Limit checks table
Replace-Value "11.72" "11.62"
Replace-Value "23.45" "23.25"
Replace-Value "1.75" "1.67"

# Shear reinforcement strength table
Replace-Value "46.9" "46.5"
Replace-Value "1.75" "1.67"
Replace-Value "11.39" "10.97"
Replace-Value "39.78" "41.41"

# Shear strength table
Replace-Value "938.0" "930.0"
Replace-Value "0.00107" "0.00108"
Replace-Value "0.834" "0.836"
Replace-Value "59.8" "59.29"
Replace-Value "291.95" "289.46"
Replace-Value "99.58" "100.7"
Replace-Value "2.21" "2.18"
